# New PO forecast model
# - Append new weeks to "Weekly Quantity"
# - Append new month to "Monthly Trend"
# - Add a new "PO Forecast" sheet with the forecast series

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet 1: "Weekly Quantity" - append two new weekly rows
# ---------------------------------------------------------------------
$wsWeekly = $wb.Worksheets.Item(1)

# Row 3
$wsWeekly.Range("A2").Copy()
$wsWeekly.Range("A3").PasteSpecial(-4122)
$wsWeekly.Range("A3").Value = 45669.99999999999
$wsWeekly.Range("B3").Value = 4

# Row 4
$wsWeekly.Range("A2").Copy()
$wsWeekly.Range("A4").PasteSpecial(-4122)
$wsWeekly.Range("A4").Value = 45676.99999999999
$wsWeekly.Range("B4").Value = 1

# ---------------------------------------------------------------------
# Sheet 2: "Monthly Trend" - append one new monthly row
# ---------------------------------------------------------------------
$wsMonthly = $wb.Worksheets.Item(2)

$wsMonthly.Range("A2").Copy()
$wsMonthly.Range("A3").PasteSpecial(-4122)
$wsMonthly.Range("A3").Value = 45688.99999999999
$wsMonthly.Range("B3").Value = 5

# ---------------------------------------------------------------------
# Sheet 3 (new): "PO Forecast"
# ---------------------------------------------------------------------
$sheetCount = $wb.Worksheets.Count
$afterSheet = $wb.Worksheets.Item($sheetCount)
$wsForecast = $wb.Worksheets.Add($null, $afterSheet)
$wsForecast.Name = "PO Forecast"

# Match page margins / outline props of the other sheets as closely as possible
$wsForecast.PageSetup().LeftMargin = 54
$wsForecast.PageSetup().RightMargin = 54
$wsForecast.PageSetup().TopMargin = 72
$wsForecast.PageSetup().BottomMargin = 72
$wsForecast.PageSetup().HeaderMargin = 36
$wsForecast.PageSetup().FooterMargin = 36
$wsForecast.Outline().SummaryRow = 1
$wsForecast.Outline().SummaryColumn = 1

# Header row (copy the bold/bordered/centered header style from sheet1)
$wsWeekly.Range("A1:B1").Copy()
$wsForecast.Range("A1:B1").PasteSpecial(-4122)
$wsForecast.Range("A1").Value = "ds"
$wsForecast.Range("B1").Value = "PO_Forecast"

# Data rows - reuse the date-styled cell (A2 on sheet1) as the style source
$dates = @(45613.99999999999, 45669.99999999999, 45676.99999999999, 45683.99999999999, 45690.99999999999, 45697.99999999999, 45704.99999999999, 45711.99999999999, 45718.99999999999, 45725.99999999999, 45732.99999999999)
$values = @(38, 5, 0, 0, 0, 0, 0, 0, 0, 0, 0)

for ($i = 0; $i -lt $dates.Length; $i++) {
    $row = $i + 2
    $wsWeekly.Range("A2").Copy()
    $wsForecast.Range("A$row").PasteSpecial(-4122)
    $wsForecast.Range("A$row").Value = $dates[$i]
    $wsForecast.Range("B$row").Value = $values[$i]
}
